$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 3972
$ws.Cells.Item(34, 9).Value = 3972
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 3972
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -3769
$ws.Cells.Item(36, 8).Value = 3972
$ws.Cells.Item(36, 9).Value = 3972
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 3972
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -3257
$ws.Cells.Item(80, 8).Value = 660.63635
$ws.Cells.Item(80, 9).Value = 727
$ws.Cells.Item(80, 10).Value = 605.3333
$ws.Cells.Item(80, 11).Value = 2181
$ws.Cells.Item(80, 12).Value = 1815.9999
$ws.Cells.Item(80, 13).Value = -1183
$ws.Cells.Item(80, 14).Value = -3811.9999
$ws.Cells.Item(83, 8).Value = 660.63635
$ws.Cells.Item(83, 9).Value = 727
$ws.Cells.Item(83, 10).Value = 605.3333
$ws.Cells.Item(83, 11).Value = 6543
$ws.Cells.Item(83, 12).Value = 5447.9997
$ws.Cells.Item(83, 13).Value = -1551
$ws.Cells.Item(83, 14).Value = -15431.9997
$ws.Cells.Item(111, 8).Value = 1625.6451
$ws.Cells.Item(111, 9).Value = 505.25
$ws.Cells.Item(111, 10).Value = 1791.6296
$ws.Cells.Item(111, 11).Value = 1515.75
$ws.Cells.Item(111, 12).Value = 5374.8888
$ws.Cells.Item(111, 13).Value = 1551.25
$ws.Cells.Item(111, 14).Value = -11508.8888
$ws.Cells.Item(125, 8).Value = 2377.625
$ws.Cells.Item(125, 9).Value = 2451.8
$ws.Cells.Item(125, 10).Value = 2254
$ws.Cells.Item(125, 11).Value = 22066.2
$ws.Cells.Item(125, 12).Value = 20286
$ws.Cells.Item(125, 13).Value = -19606.2
$ws.Cells.Item(125, 14).Value = -25206
$ws.Cells.Item(132, 8).Value = 26020.182
$ws.Cells.Item(132, 9).Value = 31261.715
$ws.Cells.Item(132, 10).Value = 16847.5
$ws.Cells.Item(132, 11).Value = 93785.145
$ws.Cells.Item(132, 12).Value = 50542.5
$ws.Cells.Item(132, 13).Value = -91255.145
$ws.Cells.Item(132, 14).Value = -55602.5
$ws.Cells.Item(137, 8).Value = 6592.0713
$ws.Cells.Item(137, 9).Value = 2258.4
$ws.Cells.Item(137, 10).Value = 8999.666999999999
$ws.Cells.Item(137, 11).Value = 6775.200000000001
$ws.Cells.Item(137, 12).Value = 26999.001
$ws.Cells.Item(137, 13).Value = -4225.200000000001
$ws.Cells.Item(137, 14).Value = -32099.001
$ws.Cells.Item(138, 8).Value = 2566.6428
$ws.Cells.Item(138, 9).Value = 1227.75
$ws.Cells.Item(138, 10).Value = 10600
$ws.Cells.Item(138, 11).Value = 3683.25
$ws.Cells.Item(138, 12).Value = 31800
$ws.Cells.Item(138, 13).Value = 1456.75
$ws.Cells.Item(138, 14).Value = -42080

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3228409.2
$ws.Cells.Item(32, 9).Value = 603.3461
$ws.Cells.Item(32, 10).Value = 20013000
$ws.Cells.Item(32, 11).Value = 603.3461
$ws.Cells.Item(32, 12).Value = 20013000
$ws.Cells.Item(32, 13).Value = -316.3461
$ws.Cells.Item(32, 14).Value = -20013574
$ws.Cells.Item(44, 8).Value = 12441.353
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 12441.353
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 12441.353
$ws.Cells.Item(44, 14).Value = -13417.353
$ws.Cells.Item(45, 8).Value = 2486.4614
$ws.Cells.Item(45, 9).Value = 1991.5
$ws.Cells.Item(45, 10).Value = 2910.7144
$ws.Cells.Item(45, 11).Value = 1991.5
$ws.Cells.Item(45, 12).Value = 2910.7144
$ws.Cells.Item(45, 13).Value = -1614.5
$ws.Cells.Item(45, 14).Value = -3664.7144
$ws.Cells.Item(61, 8).Value = 1974.25
$ws.Cells.Item(61, 9).Value = 1863.4286
$ws.Cells.Item(61, 10).Value = 2750
$ws.Cells.Item(61, 11).Value = 1863.4286
$ws.Cells.Item(61, 12).Value = 2750
$ws.Cells.Item(61, 13).Value = -1651.4286
$ws.Cells.Item(61, 14).Value = -3174
$ws.Cells.Item(88, 8).Value = 835.4545000000001
$ws.Cells.Item(88, 9).Value = 383.16666
$ws.Cells.Item(88, 10).Value = 1378.2
$ws.Cells.Item(88, 11).Value = 383.16666
$ws.Cells.Item(88, 12).Value = 1378.2
$ws.Cells.Item(88, 13).Value = 22.83334000000002
$ws.Cells.Item(88, 14).Value = -2190.2
$ws.Cells.Item(91, 8).Value = 835.4545000000001
$ws.Cells.Item(91, 9).Value = 383.16666
$ws.Cells.Item(91, 10).Value = 1378.2
$ws.Cells.Item(91, 11).Value = 383.16666
$ws.Cells.Item(91, 12).Value = 1378.2
$ws.Cells.Item(91, 13).Value = 1020.83334
$ws.Cells.Item(91, 14).Value = -4186.2
$ws.Cells.Item(136, 8).Value = 1974.25
$ws.Cells.Item(136, 9).Value = 1863.4286
$ws.Cells.Item(136, 10).Value = 2750
$ws.Cells.Item(136, 11).Value = 5590.2858
$ws.Cells.Item(136, 12).Value = 8250
$ws.Cells.Item(136, 13).Value = -3040.2858
$ws.Cells.Item(136, 14).Value = -13350

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 750
$ws.Cells.Item(20, 9).Value = 500
$ws.Cells.Item(20, 10).Value = 1000
$ws.Cells.Item(20, 11).Value = 500
$ws.Cells.Item(20, 12).Value = 1000
$ws.Cells.Item(20, 13).Value = -253
$ws.Cells.Item(20, 14).Value = -1494
$ws.Cells.Item(105, 8).Value = 6994441
$ws.Cells.Item(105, 9).Value = 11364892
$ws.Cells.Item(105, 10).Value = 1720
$ws.Cells.Item(105, 11).Value = 11364892
$ws.Cells.Item(105, 12).Value = 1720
$ws.Cells.Item(105, 13).Value = -11363145
$ws.Cells.Item(105, 14).Value = -5214

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1306.0667
$ws.Cells.Item(22, 9).Value = 771.7778
$ws.Cells.Item(22, 10).Value = 2107.5
$ws.Cells.Item(22, 11).Value = 771.7778
$ws.Cells.Item(22, 12).Value = 2107.5
$ws.Cells.Item(22, 13).Value = -421.7778
$ws.Cells.Item(22, 14).Value = -2807.5
$ws.Cells.Item(31, 8).Value = 3207.2683
$ws.Cells.Item(31, 9).Value = 1880.8572
$ws.Cells.Item(31, 10).Value = 4600
$ws.Cells.Item(31, 11).Value = 1880.8572
$ws.Cells.Item(31, 12).Value = 4600
$ws.Cells.Item(31, 13).Value = -1585.8572
$ws.Cells.Item(31, 14).Value = -5190
$ws.Cells.Item(34, 8).Value = 3207.2683
$ws.Cells.Item(34, 9).Value = 1880.8572
$ws.Cells.Item(34, 10).Value = 4600
$ws.Cells.Item(34, 11).Value = 1880.8572
$ws.Cells.Item(34, 12).Value = 4600
$ws.Cells.Item(34, 13).Value = -1678.8572
$ws.Cells.Item(34, 14).Value = -5004
$ws.Cells.Item(58, 8).Value = 2847.2917
$ws.Cells.Item(58, 9).Value = 1667.75
$ws.Cells.Item(58, 10).Value = 8745
$ws.Cells.Item(58, 11).Value = 1667.75
$ws.Cells.Item(58, 12).Value = 8745
$ws.Cells.Item(58, 13).Value = -1464.75
$ws.Cells.Item(58, 14).Value = -9151
$ws.Cells.Item(59, 8).Value = 38441.25
$ws.Cells.Item(59, 9).Value = 7000
$ws.Cells.Item(59, 10).Value = 48921.668
$ws.Cells.Item(59, 11).Value = 7000
$ws.Cells.Item(59, 12).Value = 48921.668
$ws.Cells.Item(59, 13).Value = -5855
$ws.Cells.Item(59, 14).Value = -51211.668
$ws.Cells.Item(136, 8).Value = 2847.2917
$ws.Cells.Item(136, 9).Value = 1667.75
$ws.Cells.Item(136, 10).Value = 8745
$ws.Cells.Item(136, 11).Value = 5003.25
$ws.Cells.Item(136, 12).Value = 26235
$ws.Cells.Item(136, 13).Value = -2453.25
$ws.Cells.Item(136, 14).Value = -31335

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1926082.8
$ws.Cells.Item(4, 9).Value = 2101172
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 6303516
$ws.Cells.Item(4, 12).Value = 300
$ws.Cells.Item(4, 13).Value = -6303404
$ws.Cells.Item(4, 14).Value = -524
$ws.Cells.Item(80, 8).Value = 4328.6523
$ws.Cells.Item(80, 9).Value = 4013.2222
$ws.Cells.Item(80, 10).Value = 5464.2
$ws.Cells.Item(80, 11).Value = 12039.6666
$ws.Cells.Item(80, 12).Value = 16392.6
$ws.Cells.Item(80, 13).Value = -11103.6666
$ws.Cells.Item(80, 14).Value = -18264.6
$ws.Cells.Item(83, 8).Value = 4328.6523
$ws.Cells.Item(83, 9).Value = 4013.2222
$ws.Cells.Item(83, 10).Value = 5464.2
$ws.Cells.Item(83, 11).Value = 36118.99980000001
$ws.Cells.Item(83, 12).Value = 49177.8
$ws.Cells.Item(83, 13).Value = -31438.99980000001
$ws.Cells.Item(83, 14).Value = -58537.8
$ws.Cells.Item(92, 8).Value = 5925
$ws.Cells.Item(92, 9).Value = 2000
$ws.Cells.Item(92, 10).Value = 9850
$ws.Cells.Item(92, 11).Value = 6000
$ws.Cells.Item(92, 12).Value = 29550
$ws.Cells.Item(92, 13).Value = -4752
$ws.Cells.Item(92, 14).Value = -32046
$ws.Cells.Item(140, 8).Value = 3898
$ws.Cells.Item(140, 9).Value = 3301.25
$ws.Cells.Item(140, 10).Value = 5091.5
$ws.Cells.Item(140, 11).Value = 9903.75
$ws.Cells.Item(140, 12).Value = 15274.5
$ws.Cells.Item(140, 13).Value = -4723.75
$ws.Cells.Item(140, 14).Value = -25634.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 90.57895000000001
$ws.Cells.Item(2, 9).Value = 96.833336
$ws.Cells.Item(2, 10).Value = 79.85714
$ws.Cells.Item(2, 11).Value = 96.833336
$ws.Cells.Item(2, 12).Value = 79.85714
$ws.Cells.Item(2, 13).Value = 16.166664
$ws.Cells.Item(2, 14).Value = -305.85714
$ws.Cells.Item(80, 8).Value = 8166.3335
$ws.Cells.Item(80, 9).Value = 7500
$ws.Cells.Item(80, 10).Value = 8499.5
$ws.Cells.Item(80, 11).Value = 7500
$ws.Cells.Item(80, 12).Value = 8499.5
$ws.Cells.Item(80, 13).Value = -6502
$ws.Cells.Item(80, 14).Value = -10495.5
$ws.Cells.Item(83, 8).Value = 8166.3335
$ws.Cells.Item(83, 9).Value = 7500
$ws.Cells.Item(83, 10).Value = 8499.5
$ws.Cells.Item(83, 11).Value = 37500
$ws.Cells.Item(83, 12).Value = 42497.5
$ws.Cells.Item(83, 13).Value = -32508
$ws.Cells.Item(83, 14).Value = -52481.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2477.7778
$ws.Cells.Item(22, 9).Value = 2100
$ws.Cells.Item(22, 10).Value = 2666.6667
$ws.Cells.Item(22, 11).Value = 2100
$ws.Cells.Item(22, 12).Value = 2666.6667
$ws.Cells.Item(22, 13).Value = -1805
$ws.Cells.Item(22, 14).Value = -3256.6667
$ws.Cells.Item(27, 8).Value = 2477.7778
$ws.Cells.Item(27, 9).Value = 2100
$ws.Cells.Item(27, 10).Value = 2666.6667
$ws.Cells.Item(27, 11).Value = 2100
$ws.Cells.Item(27, 12).Value = 2666.6667
$ws.Cells.Item(27, 13).Value = -1993
$ws.Cells.Item(27, 14).Value = -2880.6667
$ws.Cells.Item(55, 8).Value = 1076.6666
$ws.Cells.Item(55, 9).Value = 1076.6666
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = 1076.6666
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 13).Value = -903.6666
$ws.Cells.Item(55, 14).ClearContents()
$ws.Cells.Item(68, 8).Value = 8333.556
$ws.Cells.Item(68, 9).Value = 6000.6665
$ws.Cells.Item(68, 10).Value = 9500
$ws.Cells.Item(68, 11).Value = 6000.6665
$ws.Cells.Item(68, 12).Value = 9500
$ws.Cells.Item(68, 13).Value = -5251.6665
$ws.Cells.Item(68, 14).Value = -10998
$ws.Cells.Item(71, 8).Value = 8333.556
$ws.Cells.Item(71, 9).Value = 6000.6665
$ws.Cells.Item(71, 10).Value = 9500
$ws.Cells.Item(71, 11).Value = 30003.3325
$ws.Cells.Item(71, 12).Value = 47500
$ws.Cells.Item(71, 13).Value = -26259.3325
$ws.Cells.Item(71, 14).Value = -54988
$ws.Cells.Item(82, 8).Value = 3422.8333
$ws.Cells.Item(82, 9).Value = 1054.1666
$ws.Cells.Item(82, 10).Value = 5791.5
$ws.Cells.Item(82, 11).Value = 1054.1666
$ws.Cells.Item(82, 12).Value = 5791.5
$ws.Cells.Item(82, 13).Value = -693.1666
$ws.Cells.Item(82, 14).Value = -6513.5
$ws.Cells.Item(85, 8).Value = 3422.8333
$ws.Cells.Item(85, 9).Value = 1054.1666
$ws.Cells.Item(85, 10).Value = 5791.5
$ws.Cells.Item(85, 11).Value = 1054.1666
$ws.Cells.Item(85, 12).Value = 5791.5
$ws.Cells.Item(85, 13).Value = 193.8334
$ws.Cells.Item(85, 14).Value = -8287.5
$ws.Cells.Item(123, 8).Value = 250000
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 250000
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 250000
$ws.Cells.Item(123, 14).Value = -259800
$ws.Cells.Item(136, 8).Value = 2058.6667
$ws.Cells.Item(136, 9).Value = 1860.4
$ws.Cells.Item(136, 10).Value = 3050
$ws.Cells.Item(136, 11).Value = 5581.200000000001
$ws.Cells.Item(136, 12).Value = 9150
$ws.Cells.Item(136, 13).Value = -3031.200000000001
$ws.Cells.Item(136, 14).Value = -14250

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 2183.4783
$ws.Cells.Item(136, 9).Value = 1851.6666
$ws.Cells.Item(136, 10).Value = 3378
$ws.Cells.Item(136, 11).Value = 5554.9998
$ws.Cells.Item(136, 12).Value = 10134
$ws.Cells.Item(136, 13).Value = -3004.9998
$ws.Cells.Item(136, 14).Value = -15234
